# Week 13 logging update
# Applies cell-value updates across the YDS, OFF, DEF, ST, TURNS and PEN
# sheets, appending the week's new play-by-play numbers to the running
# space-separated logs and refreshing the derived summary totals.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# YDS sheet: append newly logged play yardages to the running lists
# ---------------------------------------------------------------------
$ydsWs = $wb.Worksheets.Item("YDS")

$ydsWs.Range("B2").Value = $ydsWs.Range("B2").Value2 + " 0 2 2 4 11 2 0 6 1 6 2 3 0 2 3 4 1 3 4 4 2 3 2 2"
$ydsWs.Range("C2").Value = $ydsWs.Range("C2").Value2 + " 3 3 10 -2 0 1 0 7 0 23 17 1 0 9 13 3"
$ydsWs.Range("B3").Value = $ydsWs.Range("B3").Value2 + " 9 5 11 1 4 12 8 6 -2 6 8 4 8 5 9 2 9 17 4 6 5 3 17 25 11 13 3 2 16 17"
$ydsWs.Range("C3").Value = $ydsWs.Range("C3").Value2 + " 3 15 -3 7 -5 4 18 20 -1 5 12 16 6 11 17 8 3 9 8 2 10 11 11"

# ---------------------------------------------------------------------
# OFF sheet: updated running totals
# ---------------------------------------------------------------------
$offWs = $wb.Worksheets.Item("OFF")

$offWs.Range("C2").Value = 141
$offWs.Range("E2").Value = 11
$offWs.Range("F2").Value = 50
$offWs.Range("G2").Value = 39
$offWs.Range("J2").Value = 21
$offWs.Range("L2").Value = 266
$offWs.Range("M2").Value = 179
$offWs.Range("Q2").Value = 461

$offWs.Range("C3").Value = 151
$offWs.Range("E3").Value = 23
$offWs.Range("F3").Value = 91
$offWs.Range("G3").Value = 38
$offWs.Range("H3").Value = 27
$offWs.Range("I3").Value = 45
$offWs.Range("J3").Value = 60
$offWs.Range("N3").Value = 22

# ---------------------------------------------------------------------
# DEF sheet: updated running totals
# ---------------------------------------------------------------------
$defWs = $wb.Worksheets.Item("DEF")

$defWs.Range("C2").Value = 146
$defWs.Range("F2").Value = 48
$defWs.Range("G2").Value = 24
$defWs.Range("H2").Value = 2
$defWs.Range("I2").Value = 5
$defWs.Range("L2").Value = 261
$defWs.Range("M2").Value = 148
$defWs.Range("O2").Value = 13
$defWs.Range("P2").Value = 9
$defWs.Range("Q2").Value = 437

$defWs.Range("C3").Value = 163
$defWs.Range("D3").Value = 3
$defWs.Range("E3").Value = 40
$defWs.Range("F3").Value = 101
$defWs.Range("H3").Value = 34
$defWs.Range("I3").Value = 58
$defWs.Range("J3").Value = 41
$defWs.Range("N3").Value = 21

# ---------------------------------------------------------------------
# ST sheet: updated running totals and appended kick-distance logs
# ---------------------------------------------------------------------
$stWs = $wb.Worksheets.Item("ST")

$stWs.Range("B2").Value = 59
$stWs.Range("D2").Value = 60
$stWs.Range("F2").Value = 128
$stWs.Range("G2").Value = 125
$stWs.Range("L2").Value = 40
$stWs.Range("M2").Value = 30
$stWs.Range("N2").Value = 19

$stWs.Range("B3").Value = 40

$stWs.Range("D3").Value = $stWs.Range("D3").Value2 + " 53 50 65 39 37 46"
$stWs.Range("D4").Value = $stWs.Range("D4").Value2 + " 0 15 11 0 -4 0"
$stWs.Range("D5").Value = $stWs.Range("D5").Value2 + " 0 0 7 0 0 0"
$stWs.Range("B6").Value = $stWs.Range("B6").Value2 + " 15 17"

# ---------------------------------------------------------------------
# TURNS sheet: updated running totals
# ---------------------------------------------------------------------
$turnsWs = $wb.Worksheets.Item("TURNS")

$turnsWs.Range("C2").Value = 10
$turnsWs.Range("D3").Value = 8

# ---------------------------------------------------------------------
# PEN sheet: updated running totals
# ---------------------------------------------------------------------
$penWs = $wb.Worksheets.Item("PEN")

$penWs.Range("B2").Value = 19
$penWs.Range("D2").Value = 5

Write-Host "Week 13 logging applied"
